$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the sensor columns (C:H) for rows 2-21 down by one row, taking the
# value that used to be in the row above (i.e. row r gets what was in row r-1).
# We copy from the bottom upward is not required since Copy/Paste captures the
# whole source range at once before writing to destination.
$src = $ws.Range("C2:H21")
$dst = $ws.Range("C3:H22")
$src.Copy($dst)

# Write the brand-new first data row (timestamp 0, label "falling") values.
$ws.Range("C2").Value = -3.729709470272064
$ws.Range("D2").Value = 9.457800364494323
$ws.Range("E2").Value = 0.187229474633932
$ws.Range("F2").Value = 0.0209221355617046
$ws.Range("G2").Value = -0.0198531206697225
$ws.Range("H2").Value = -0.0239764600992202

# Remove the now-duplicated last row (row 22), shifting nothing else because
# it is the last row of data.
$ws.Rows("22").Delete()

$excel.CutCopyMode = $false
